$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 497, pushing the existing
# rows 497-504 down to 499-506 (matches the diff's row-shift pattern).
$ws.Rows.Item(497).Insert()
$ws.Rows.Item(497).Insert()

# Copy the date number format (style used by column D, e.g. D496) onto
# the two new D cells so the new dates render the same way as the rest
# of the column.
$ws.Range("D496").Copy()
$ws.Range("D497:D498").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 497: Lechuga / Escarola / Segunda
$ws.Cells.Item(497, 1).Value = 4
$ws.Cells.Item(497, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(497, 3).Value = "Los Lagos"
$ws.Cells.Item(497, 4).Value = 44656
$ws.Cells.Item(497, 5).Value = 10
$ws.Cells.Item(497, 6).Value = 100112033
$ws.Cells.Item(497, 7).Value = "Lechuga"
$ws.Cells.Item(497, 8).Value = "Escarola"
$ws.Cells.Item(497, 9).Value = "Segunda"
$ws.Cells.Item(497, 10).Value = 300
$ws.Cells.Item(497, 11).Value = 11000
$ws.Cells.Item(497, 12).Value = 11000
$ws.Cells.Item(497, 13).Value = 11000
$ws.Cells.Item(497, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(497, 15).Value = "Región Metropolitana"
$ws.Cells.Item(497, 16).Value = 611
$ws.Cells.Item(497, 17).Value = 18
$ws.Cells.Item(497, 18).Value = "Hortaliza"

# New row 498: Lechuga / Escarola / Segunda
$ws.Cells.Item(498, 1).Value = 4
$ws.Cells.Item(498, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(498, 3).Value = "Los Lagos"
$ws.Cells.Item(498, 4).Value = 44656
$ws.Cells.Item(498, 5).Value = 10
$ws.Cells.Item(498, 6).Value = 100112033
$ws.Cells.Item(498, 7).Value = "Lechuga"
$ws.Cells.Item(498, 8).Value = "Escarola"
$ws.Cells.Item(498, 9).Value = "Segunda"
$ws.Cells.Item(498, 10).Value = 200
$ws.Cells.Item(498, 11).Value = 9000
$ws.Cells.Item(498, 12).Value = 9000
$ws.Cells.Item(498, 13).Value = 9000
$ws.Cells.Item(498, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(498, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(498, 16).Value = 500
$ws.Cells.Item(498, 17).Value = 18
$ws.Cells.Item(498, 18).Value = "Hortaliza"
